# Update MCDA default weights and ranges
#
# Sheet "3-state" (sheet1): replace the old "Aes combined" row (row 8) with
# "Elevated alanine transaminase" and append 9 more adverse-event criteria
# rows (rows 9-17), each with weight formula =5/10 and new min values.
#
# Sheet "4-state" (sheet2): same, but the AE rows start one row lower
# (row 9 replaced, rows 10-18 appended) because this sheet has the extra
# "PFS with 2L treatment" criterion.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("3-state")
$ws2 = $wb.Worksheets.Item("4-state")

# New adverse-event criteria, in the order they should first appear
# (this also controls shared-string insertion order).
$aeNames = @(
    "Elevated alanine transaminase",
    "Elevated aspartate transaminase",
    "Diarrhea",
    "Dry skin",
    "Eye problems",
    "Paronychia",
    "Pneumonitis",
    "Pruritis",
    "Rash",
    "Stomatitis"
)
$aeMins = @(0.26, 0.2, 0.3, 0.03, 0.02, 0.1, 0.06, 0.015, 0.65, 0.14)
$aeMaxs = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)

# ---------------------------------------------------------------------
# Sheet "3-state": old data rows are 2..8 (row 8 = "Aes combined").
# New layout: rows 2..7 stay as before, rows 8..17 become the AE block.
# ---------------------------------------------------------------------
$startRow1 = 8
for ($i = 0; $i -lt $aeNames.Count; $i++) {
    $r = $startRow1 + $i
    $ws1.Cells.Item($r, 1).Value = $aeNames[$i]
    $ws1.Cells.Item($r, 3).Value = $aeMins[$i]
    $ws1.Cells.Item($r, 4).Value = $aeMaxs[$i]
}
# First AE weight cell gets its own (non-shared) formula ...
$ws1.Range("B8").Formula = "=5/10"
# ... the remaining AE weight cells share one formula definition.
$ws1.Range("B9:B17").Formula = "=5/10"
# Apply the 2-decimal number format only to the newly written weight cells.
$ws1.Range("B8:B17").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# Sheet "4-state": old data rows are 2..9 (row 9 = "Aes combined").
# New layout: rows 2..8 stay as before, rows 9..18 become the AE block.
# ---------------------------------------------------------------------
$startRow2 = 9
for ($i = 0; $i -lt $aeNames.Count; $i++) {
    $r = $startRow2 + $i
    $ws2.Cells.Item($r, 1).Value = $aeNames[$i]
    $ws2.Cells.Item($r, 3).Value = $aeMins[$i]
    $ws2.Cells.Item($r, 4).Value = $aeMaxs[$i]
}
$ws2.Range("B9").Formula = "=5/10"
$ws2.Range("B10:B18").Formula = "=5/10"
# On this sheet the whole weight column (including the pre-existing rows)
# picks up the 2-decimal format.
$ws2.Range("B1:B18").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# View / selection bookkeeping to match the edited workbook: sheet
# "3-state" becomes the active tab, "4-state" loses it.
# ---------------------------------------------------------------------
$ws2.Range("A9:D18").Select()
$ws1.Activate()
$ws1.Range("B4").Select()
